# Rename "Gamepad controller" -> "Input controller" wherever it appears,
# preserving the existing run/formatting structure (font size, fill, etc.).
$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        $hasTF = $sh.HasTextFrame
        if ($hasTF -eq -1 -or $hasTF -eq $true) {
            $tr = $sh.TextFrame.TextRange
            $txt = $tr.Text
            if ($txt -eq "Gamepad controller") {
                $len = $tr.Length
                $splitAt = 7   # length of the word "Gamepad"

                $firstPart = $tr.Characters(1, $splitAt)
                $secondPart = $tr.Characters($splitAt + 1, $len - $splitAt)

                # Only treat this as two independent runs if they really do
                # carry different character formatting (e.g. differing font
                # sizes); otherwise it is a single run and must stay that way.
                $isTwoRuns = ($firstPart.Font.Size -ne $secondPart.Font.Size)

                if ($isTwoRuns) {
                    $firstPart.Text = "Input "
                    $newLen = $tr.Length
                    $secondPart2 = $tr.Characters(7, $newLen - 6)
                    $secondPart2.Text = "controller"
                } else {
                    $tr.Text = "Input controller"
                }
            }
        }
    }
}
